$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1
$ws.Range("B3").Value = 3.7
$ws.Range("C3").Value = 7.1
$ws.Range("C4").Value = 8.699999999999999
$ws.Range("C5").Value = 12.5
$ws.Range("C6").Value = 7.5
$ws.Range("C7").Value = 0.6
$ws.Range("C13").Value = 5.2
